$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4206.353
$ws.Range("J17").Value = 4206.353
$ws.Range("L17").Value = 12619.059
$ws.Range("N17").Value = -12955.059

$ws.Range("H53").Value = 232.33333
$ws.Range("I53").Value = 149.66667
$ws.Range("K53").Value = 149.66667
$ws.Range("M53").Value = 487.33333

$ws.Range("H62").Value = 7640.625
$ws.Range("I62").Value = 6866.5
$ws.Range("J62").Value = 7898.6665
$ws.Range("K62").Value = 6866.5
$ws.Range("L62").Value = 7898.6665
$ws.Range("M62").Value = -6242.5
$ws.Range("N62").Value = -9146.666499999999

$ws.Range("H65").Value = 7640.625
$ws.Range("I65").Value = 6866.5
$ws.Range("J65").Value = 7898.6665
$ws.Range("K65").Value = 34332.5
$ws.Range("L65").Value = 39493.3325
$ws.Range("M65").Value = -31212.5
$ws.Range("N65").Value = -45733.3325

$ws.Range("H111").Value = 1564.8182
$ws.Range("I111").Value = 626
$ws.Range("J111").Value = 2691.4
$ws.Range("K111").Value = 1878
$ws.Range("L111").Value = 8074.200000000001
$ws.Range("M111").Value = 1189
$ws.Range("N111").Value = -14208.2

$ws.Range("H116").Value = 9299
$ws.Range("I116").Value = 11699
$ws.Range("K116").Value = 11699
$ws.Range("M116").Value = -8257

$ws.Range("H132").Value = 755.65
$ws.Range("I132").Value = 724.55554
$ws.Range("K132").Value = 2173.66662
$ws.Range("M132").Value = 356.33338

$ws.Range("H137").Value = 1601.7142
$ws.Range("I137").Value = 1609.5385
$ws.Range("J137").Value = 1500
$ws.Range("K137").Value = 4828.6155
$ws.Range("L137").Value = 4500
$ws.Range("M137").Value = -2278.6155
$ws.Range("N137").Value = -9600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11810.909
$ws.Range("I32").Value = 11992.5
$ws.Range("J32").Value = 9995
$ws.Range("K32").Value = 11992.5
$ws.Range("L32").Value = 9995
$ws.Range("M32").Value = -11705.5
$ws.Range("N32").Value = -10569

$ws.Range("H36").Value = 4912.3335
$ws.Range("I36").Value = 4912.3335
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 4912.3335
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -4566.3335
$ws.Range("N36").ClearContents()

$ws.Range("H45").Value = 4390.727
$ws.Range("I45").Value = 3474.4443
$ws.Range("K45").Value = 3474.4443
$ws.Range("M45").Value = -3097.4443

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 397
$ws.Range("J80").Value = 224.45454
$ws.Range("L80").Value = 224.45454
$ws.Range("N80").Value = -2220.45454

$ws.Range("H83").Value = 397
$ws.Range("J83").Value = 224.45454
$ws.Range("L83").Value = 1122.2727
$ws.Range("N83").Value = -11106.2727

$ws.Range("H86").Value = 3076.375
$ws.Range("I86").Value = 657.7143
$ws.Range("K86").Value = 657.7143
$ws.Range("M86").Value = 465.2857

$ws.Range("H89").Value = 3076.375
$ws.Range("I89").Value = 657.7143
$ws.Range("K89").Value = 3288.5715
$ws.Range("M89").Value = 2327.4285

$ws.Range("H133").Value = 58748.75
$ws.Range("I133").Value = 49998.5
$ws.Range("K133").Value = 49998.5
$ws.Range("M133").Value = -44938.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2233.7693
$ws.Range("I31").Value = 1859.2858
$ws.Range("J31").Value = 2670.6667
$ws.Range("K31").Value = 1859.2858
$ws.Range("L31").Value = 2670.6667
$ws.Range("M31").Value = -1564.2858
$ws.Range("N31").Value = -3260.6667

$ws.Range("H34").Value = 2233.7693
$ws.Range("I34").Value = 1859.2858
$ws.Range("J34").Value = 2670.6667
$ws.Range("K34").Value = 1859.2858
$ws.Range("L34").Value = 2670.6667
$ws.Range("M34").Value = -1657.2858
$ws.Range("N34").Value = -3074.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8333437
$ws.Range("I4").Value = 8333437
$ws.Range("K4").Value = 25000311
$ws.Range("M4").Value = -25000199

$ws.Range("H5").Value = 5398.5
$ws.Range("I5").Value = 5398
$ws.Range("K5").Value = 16194
$ws.Range("M5").Value = -16082

$ws.Range("H8").Value = 316.5
$ws.Range("I8").Value = 316.5
$ws.Range("K8").Value = 949.5
$ws.Range("M8").Value = -810.5

$ws.Range("H11").Value = 31272330
$ws.Range("I11").Value = 34115196
$ws.Range("J11").Value = 800
$ws.Range("K11").Value = 102345588
$ws.Range("L11").Value = 2400
$ws.Range("M11").Value = -102345448
$ws.Range("N11").Value = -2680

$ws.Range("H23").Value = 60000064
$ws.Range("J23").Value = 41.5
$ws.Range("L23").Value = 124.5
$ws.Range("N23").Value = -594.5

$ws.Range("H26").Value = 898
$ws.Range("I26").Value = 898
$ws.Range("K26").Value = 2694
$ws.Range("M26").Value = -2406

$ws.Range("H37").Value = 79975
$ws.Range("J37").Value = 79975
$ws.Range("L37").Value = 239925
$ws.Range("N37").Value = -240149

$ws.Range("H64").Value = 18890
$ws.Range("J64").Value = 18890
$ws.Range("L64").Value = 56670
$ws.Range("N64").Value = -57210

$ws.Range("H67").Value = 18890
$ws.Range("J67").Value = 18890
$ws.Range("L67").Value = 56670
$ws.Range("N67").Value = -58542

$ws.Range("H68").Value = 16458.125
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 16458.125
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 49374.375
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -50996.375

$ws.Range("H71").Value = 16458.125
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 16458.125
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 148123.125
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -156235.125

$ws.Range("H86").Value = 6565.75
$ws.Range("J86").Value = 6565.75
$ws.Range("L86").Value = 19697.25
$ws.Range("N86").Value = -22069.25

$ws.Range("H89").Value = 6565.75
$ws.Range("J89").Value = 6565.75
$ws.Range("L89").Value = 59091.75
$ws.Range("N89").Value = -70947.75

$ws.Range("H135").Value = 5398.5
$ws.Range("I135").Value = 5398
$ws.Range("K135").Value = 48582
$ws.Range("M135").Value = -46047

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 48000
$ws.Range("J34").Value = 48000
$ws.Range("L34").Value = 48000
$ws.Range("N34").Value = -48536

$ws.Range("H70").Value = 7904.8335
$ws.Range("I70").Value = 7772.4546
$ws.Range("J70").Value = 8112.857
$ws.Range("K70").Value = 7772.4546
$ws.Range("L70").Value = 8112.857
$ws.Range("M70").Value = -7502.4546
$ws.Range("N70").Value = -8652.857

$ws.Range("H73").Value = 7904.8335
$ws.Range("I73").Value = 7772.4546
$ws.Range("J73").Value = 8112.857
$ws.Range("K73").Value = 7772.4546
$ws.Range("L73").Value = 8112.857
$ws.Range("M73").Value = -6836.4546
$ws.Range("N73").Value = -9984.857

$ws.Range("H76").Value = 48000
$ws.Range("J76").Value = 48000
$ws.Range("L76").Value = 48000
$ws.Range("N76").Value = -48630

$ws.Range("H79").Value = 48000
$ws.Range("J79").Value = 48000
$ws.Range("L79").Value = 48000
$ws.Range("N79").Value = -50184

$ws.Range("H126").Value = 3920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1787
$ws.Range("I22").Value = 1750
$ws.Range("K22").Value = 1750
$ws.Range("M22").Value = -1455

$ws.Range("H27").Value = 1787
$ws.Range("I27").Value = 1750
$ws.Range("K27").Value = 1750
$ws.Range("M27").Value = -1643

$ws.Range("H40").Value = 4201.25
$ws.Range("I40").Value = 2410.5
$ws.Range("K40").Value = 2410.5
$ws.Range("M40").Value = -2274.5

$ws.Range("H55").Value = 313.13333
$ws.Range("I55").Value = 224.66667
$ws.Range("J55").Value = 667
$ws.Range("K55").Value = 224.66667
$ws.Range("L55").Value = 667
$ws.Range("M55").Value = -51.66667000000001
$ws.Range("N55").Value = -1013

$ws.Range("H122").Value = 7336.815
$ws.Range("I122").Value = 8692.462
$ws.Range("K122").Value = 26077.386
$ws.Range("M122").Value = -23627.386

$ws.Range("H136").Value = 3007.5
$ws.Range("I136").Value = 2394.6924
$ws.Range("J136").Value = 4600.8
$ws.Range("K136").Value = 7184.0772
$ws.Range("L136").Value = 13802.4
$ws.Range("M136").Value = -4634.0772
$ws.Range("N136").Value = -18902.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 562499.5
$ws.Range("J2").Value = 83332.664
$ws.Range("L2").Value = 83332.664
$ws.Range("N2").Value = -83556.664

$ws.Range("H62").Value = 18587.273
$ws.Range("J62").Value = 14837.286
$ws.Range("L62").Value = 14837.286
$ws.Range("N62").Value = -16085.286

$ws.Range("H65").Value = 18587.273
$ws.Range("J65").Value = 14837.286
$ws.Range("L65").Value = 74186.42999999999
$ws.Range("N65").Value = -80426.42999999999

$ws.Range("H70").Value = 90095
$ws.Range("I70").Value = 90095
$ws.Range("K70").Value = 90095
$ws.Range("M70").Value = -89780

$ws.Range("H73").Value = 90095
$ws.Range("I73").Value = 90095
$ws.Range("K73").Value = 90095
$ws.Range("M73").Value = -89003

$ws.Range("H122").Value = 6095.5
$ws.Range("I122").Value = 831.3333
$ws.Range("J122").Value = 9254
$ws.Range("K122").Value = 2493.9999
$ws.Range("L122").Value = 27762
$ws.Range("M122").Value = -43.9998999999998
$ws.Range("N122").Value = -32662
